$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 23 data (Date, Start Time, End Time)
$ws.Range("A23").Value = 43058
$ws.Range("B23").Value = 0.84027777777777779
$ws.Range("C23").Value = 0.92361111111111116

# Update selection on the sheet view
$ws.Range("V22").Select()
